$d = $word.ActiveDocument

$replacements = @(
    @{old="854×7=5978"; new="824×2=1648"},
    @{old="779×8=6232"; new="860×5=4300"},
    @{old="188×4=752"; new="618×8=4944"},
    @{old="283×6=1698"; new="421×8=3368"},
    @{old="960×8=7680"; new="628×4=2512"},
    @{old="823×4=3292"; new="708×7=4956"},
    @{old="489×2=978"; new="620×2=1240"},
    @{old="604×3=1812"; new="961×4=3844"},
    @{old="721×4=2884"; new="116×6=696"},
    @{old="628×3=1884"; new="971×2=1942"},
    @{old="581×5=2905"; new="516×5=2580"},
    @{old="539×4=2156"; new="625×8=5000"},
    @{old="643×4=2572"; new="278×8=2224"},
    @{old="867×5=4335"; new="520×3=1560"},
    @{old="746×8=5968"; new="752×6=4512"},
    @{old="690×8=5520"; new="829×3=2487"},
    @{old="552×4=2208"; new="823×9=7407"},
    @{old="495×9=4455"; new="911×3=2733"},
    @{old="482×5=2410"; new="568×3=1704"},
    @{old="807×9=7263"; new="984×8=7872"},
    @{old="523×3=1569"; new="153×8=1224"},
    @{old="131×8=1048"; new="771×4=3084"},
    @{old="182×7=1274"; new="273×5=1365"},
    @{old="630×8=5040"; new="523×6=3138"},
    @{old="492×5=2460"; new="875×8=7000"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
